$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-30 and add new row 31 with updated weekly price data
$data = @(
    @(2, 44978, 1000, 1800, 2000, 1900, 633),
    @(3, 45006, 1100, 2000, 2500, 2250, 750),
    @(4, 45007, 1160, 2000, 2500, 2250, 750),
    @(5, 44911, 700, 1800, 2000, 1900, 633),
    @(6, 45035, 1100, 2000, 2500, 2250, 750),
    @(7, 44964, 1000, 2000, 2500, 2250, 750),
    @(8, 44965, 1120, 2000, 2500, 2250, 750),
    @(9, 45070, 800, 2000, 2500, 2250, 750),
    @(10, 45091, 800, 2000, 2500, 2250, 750),
    @(11, 45013, 1100, 2000, 2500, 2250, 750),
    @(12, 44910, 1000, 1800, 2000, 1900, 633),
    @(13, 45077, 760, 2000, 2500, 2250, 750),
    @(14, 45062, 1100, 2000, 2500, 2250, 750),
    @(15, 44999, 1100, 2000, 2500, 2250, 750),
    @(16, 44881, 500, 1900, 2000, 1950, 650),
    @(17, 44985, 1000, 2000, 2500, 2250, 750),
    @(18, 44883, 500, 1800, 2000, 1900, 633),
    @(19, 44992, 1040, 2000, 2500, 2250, 750),
    @(20, 45020, 1200, 2000, 2500, 2250, 750),
    @(21, 44970, 800, 2000, 2500, 2250, 750),
    @(22, 44951, 800, 2000, 2500, 2250, 750),
    @(23, 45084, 900, 2000, 2500, 2250, 750),
    @(24, 44953, 1000, 2000, 2500, 2250, 750),
    @(25, 44685, 400, 1500, 2000, 1750, 583),
    @(26, 45028, 1000, 2000, 2500, 2250, 750),
    @(27, 44848, 1000, 1500, 2000, 1750, 583),
    @(28, 44827, 1200, 2000, 2500, 2250, 750),
    @(29, 44971, 1000, 2000, 2500, 2250, 750),
    @(30, 45041, 1160, 2000, 2500, 2250, 750),
    @(31, 45034, 1100, 2000, 2500, 2250, 750),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]  # J Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]  # K Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]  # L Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]  # M Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]  # P Precio $/Kg
}

# Copy the date number format from an existing date cell (D2) onto the new rows date cell
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# Fill the static (unchanged-across-rows) columns for the newly added row 31
$ws.Cells.Item(31, 1).Value = 8                                       # A Mercado ID
$ws.Cells.Item(31, 2).Value = "Terminal La Palmera de La Serena"      # B Mercado
$ws.Cells.Item(31, 3).Value = "Coquimbo"                              # C Region
$ws.Cells.Item(31, 5).Value = 4                                       # E Codreg
$ws.Cells.Item(31, 6).Value = 100112039                               # F Categoria ID
$ws.Cells.Item(31, 7).Value = "Ciboulette"                            # G Categoria
$ws.Cells.Item(31, 8).Value = "Sin especificar"                       # H Variedad
$ws.Cells.Item(31, 9).Value = "Primera"                               # I Calidad
$ws.Cells.Item(31, 14).Value = '$/docena de atados'                   # N Unidad de comercializacion
$ws.Cells.Item(31, 15).Value = "Provincia del Elquí"                  # O Origen
$ws.Cells.Item(31, 17).Value = 3                                      # Q Kg o Unidades
$ws.Cells.Item(31, 18).Value = "Hortaliza"                            # R Clasificacion
